# Auto-generated Excel COM-interop script to apply the scheduled market-data refresh.
# Updates hardcoded price/profit values (columns H-N) on several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 126.64286
$ws.Range("I55").Value = 97.625
$ws.Range("J55").Value = 165.33333
$ws.Range("K55").Value = 97.625
$ws.Range("L55").Value = 165.33333
$ws.Range("M55").Value = 116.375
$ws.Range("N55").Value = -593.3333299999999
$ws.Range("H62").Value = 1914.4324
$ws.Range("I62").Value = 1407.4783
$ws.Range("J62").Value = 2747.2856
$ws.Range("K62").Value = 1407.4783
$ws.Range("L62").Value = 2747.2856
$ws.Range("M62").Value = -783.4783
$ws.Range("N62").Value = -3995.2856
$ws.Range("H65").Value = 1914.4324
$ws.Range("I65").Value = 1407.4783
$ws.Range("J65").Value = 2747.2856
$ws.Range("K65").Value = 7037.3915
$ws.Range("L65").Value = 13736.428
$ws.Range("M65").Value = -3917.3915
$ws.Range("N65").Value = -19976.428
$ws.Range("H137").Value = 1853425.6
$ws.Range("I137").Value = 2381971
$ws.Range("J137").Value = 3516.8333
$ws.Range("K137").Value = 7145913
$ws.Range("L137").Value = 10550.4999
$ws.Range("M137").Value = -7143363
$ws.Range("N137").Value = -15650.4999
$ws.Range("H138").Value = 2605653.2
$ws.Range("I138").Value = 1167.3158
$ws.Range("J138").Value = 23813610
$ws.Range("K138").Value = 3501.9474
$ws.Range("L138").Value = 71440830
$ws.Range("M138").Value = 1638.0526
$ws.Range("N138").Value = -71451110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6869.31
$ws.Range("I32").Value = 3768.1392
$ws.Range("J32").Value = 18535.62
$ws.Range("K32").Value = 3768.1392
$ws.Range("L32").Value = 18535.62
$ws.Range("M32").Value = -3481.1392
$ws.Range("N32").Value = -19109.62
$ws.Range("H61").Value = 20450244
$ws.Range("I61").Value = 24415664
$ws.Range("J61").Value = 127463.5
$ws.Range("K61").Value = 24415664
$ws.Range("L61").Value = 127463.5
$ws.Range("M61").Value = -24415452
$ws.Range("N61").Value = -127887.5
$ws.Range("H74").Value = 7799694
$ws.Range("I74").Value = 10132373
$ws.Range("J74").Value = 101852.4
$ws.Range("K74").Value = 10132373
$ws.Range("L74").Value = 101852.4
$ws.Range("M74").Value = -10131499
$ws.Range("N74").Value = -103600.4
$ws.Range("H77").Value = 7799694
$ws.Range("I77").Value = 10132373
$ws.Range("J77").Value = 101852.4
$ws.Range("K77").Value = 50661865
$ws.Range("L77").Value = 509262
$ws.Range("M77").Value = -50657497
$ws.Range("N77").Value = -517998
$ws.Range("H110").Value = 294959.03
$ws.Range("I110").Value = 400548
$ws.Range("J110").Value = 1656.3334
$ws.Range("K110").Value = 400548
$ws.Range("L110").Value = 1656.3334
$ws.Range("M110").Value = -398503
$ws.Range("N110").Value = -5746.3334
$ws.Range("H122").Value = 4631720
$ws.Range("I122").Value = 2273.7334
$ws.Range("J122").Value = 12347464
$ws.Range("K122").Value = 6821.2002
$ws.Range("L122").Value = 37042392
$ws.Range("M122").Value = -4371.2002
$ws.Range("N122").Value = -37047292
$ws.Range("H125").Value = 54638.46
$ws.Range("J125").Value = 54638.46
$ws.Range("L125").Value = 54638.46
$ws.Range("N125").Value = -64478.46
$ws.Range("H132").Value = 49655.14
$ws.Range("I132").Value = 37809
$ws.Range("J132").Value = 69645.5
$ws.Range("K132").Value = 113427
$ws.Range("L132").Value = 208936.5
$ws.Range("M132").Value = -110897
$ws.Range("N132").Value = -213996.5
$ws.Range("H135").Value = 48899.8
$ws.Range("J135").Value = 48899.8
$ws.Range("L135").Value = 48899.8
$ws.Range("N135").Value = -59039.8
$ws.Range("H136").Value = 20450244
$ws.Range("I136").Value = 24415664
$ws.Range("J136").Value = 127463.5
$ws.Range("K136").Value = 73246992
$ws.Range("L136").Value = 382390.5
$ws.Range("M136").Value = -73244442
$ws.Range("N136").Value = -387490.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1156.1666
$ws.Range("I20").Value = 1152.6666
$ws.Range("J20").Value = 1159.6666
$ws.Range("K20").Value = 1152.6666
$ws.Range("L20").Value = 1159.6666
$ws.Range("M20").Value = -905.6666
$ws.Range("N20").Value = -1653.6666
$ws.Range("H99").Value = 1185.125
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 1370.25
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 1370.25
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -4366.25
$ws.Range("H132").Value = 42685
$ws.Range("I132").Value = 35000
$ws.Range("J132").Value = 45246.668
$ws.Range("K132").Value = 35000
$ws.Range("L132").Value = 45246.668
$ws.Range("M132").Value = -29940
$ws.Range("N132").Value = -55366.668
$ws.Range("H134").Value = 2624.3667
$ws.Range("I134").Value = 1368.238
$ws.Range("J134").Value = 5555.3335
$ws.Range("K134").Value = 4104.714
$ws.Range("L134").Value = 16666.0005
$ws.Range("M134").Value = -1569.714
$ws.Range("N134").Value = -21736.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2012.3422
$ws.Range("I31").Value = 1102.0294
$ws.Range("J31").Value = 9750
$ws.Range("K31").Value = 1102.0294
$ws.Range("L31").Value = 9750
$ws.Range("M31").Value = -807.0293999999999
$ws.Range("N31").Value = -10340
$ws.Range("H34").Value = 2012.3422
$ws.Range("I34").Value = 1102.0294
$ws.Range("J34").Value = 9750
$ws.Range("K34").Value = 1102.0294
$ws.Range("L34").Value = 9750
$ws.Range("M34").Value = -900.0293999999999
$ws.Range("N34").Value = -10154
$ws.Range("H58").Value = 16668106
$ws.Range("I58").Value = 21740512
$ws.Range("J58").Value = 1630.8572
$ws.Range("K58").Value = 21740512
$ws.Range("L58").Value = 1630.8572
$ws.Range("M58").Value = -21740309
$ws.Range("N58").Value = -2036.8572
$ws.Range("H62").Value = 2613.9285
$ws.Range("I62").Value = 2833.3333
$ws.Range("J62").Value = 2554.0908
$ws.Range("K62").Value = 2833.3333
$ws.Range("L62").Value = 2554.0908
$ws.Range("M62").Value = -2209.3333
$ws.Range("N62").Value = -3802.0908
$ws.Range("H65").Value = 2613.9285
$ws.Range("I65").Value = 2833.3333
$ws.Range("J65").Value = 2554.0908
$ws.Range("K65").Value = 14166.6665
$ws.Range("L65").Value = 12770.454
$ws.Range("M65").Value = -11046.6665
$ws.Range("N65").Value = -19010.454
$ws.Range("H116").Value = 49900
$ws.Range("J116").Value = 49900
$ws.Range("L116").Value = 49900
$ws.Range("M116").Value = -59078
$ws.Range("H132").Value = 33532.156
$ws.Range("I132").Value = 2425.0557
$ws.Range("J132").Value = 73527
$ws.Range("K132").Value = 7275.1671
$ws.Range("L132").Value = 220581
$ws.Range("M132").Value = -4745.1671
$ws.Range("N132").Value = -225641
$ws.Range("H134").Value = 32551.486
$ws.Range("I134").Value = 1696.4073
$ws.Range("J134").Value = 136687.38
$ws.Range("K134").Value = 5089.2219
$ws.Range("L134").Value = 410062.14
$ws.Range("M134").Value = -2554.2219
$ws.Range("N134").Value = -415132.14
$ws.Range("H136").Value = 16668106
$ws.Range("I136").Value = 21740512
$ws.Range("J136").Value = 1630.8572
$ws.Range("K136").Value = 65221536
$ws.Range("L136").Value = 4892.571599999999
$ws.Range("M136").Value = -65218986
$ws.Range("N136").Value = -9992.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 359.33334
$ws.Range("I7").Value = 138.57143
$ws.Range("J7").Value = 668.4
$ws.Range("K7").Value = 415.71429
$ws.Range("L7").Value = 2005.2
$ws.Range("M7").Value = -303.71429
$ws.Range("N7").Value = -2229.2
$ws.Range("H8").Value = 33.11111
$ws.Range("I8").Value = 33.11111
$ws.Range("K8").Value = 99.33332999999999
$ws.Range("M8").Value = 39.66667000000001
$ws.Range("H80").Value = 2080
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 1600
$ws.Range("K80").Value = 12000
$ws.Range("L80").Value = 4800
$ws.Range("M80").Value = -11064
$ws.Range("N80").Value = -6672
$ws.Range("H83").Value = 2080
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 1600
$ws.Range("K83").Value = 36000
$ws.Range("L83").Value = 14400
$ws.Range("M83").Value = -31320
$ws.Range("N83").Value = -23760
$ws.Range("H92").Value = 948.5333000000001
$ws.Range("I92").Value = 950
$ws.Range("J92").Value = 947.25
$ws.Range("K92").Value = 2850
$ws.Range("L92").Value = 2841.75
$ws.Range("M92").Value = -1602
$ws.Range("N92").Value = -5337.75
$ws.Range("H113").Value = 481.26923
$ws.Range("I113").Value = 419.66666
$ws.Range("J113").Value = 513.8823
$ws.Range("K113").Value = 1258.99998
$ws.Range("L113").Value = 1541.6469
$ws.Range("M113").Value = 911.0000199999999
$ws.Range("N113").Value = -5881.6469
$ws.Range("H140").Value = 2628.7021
$ws.Range("I140").Value = 2721.111
$ws.Range("J140").Value = 2326.2727
$ws.Range("K140").Value = 8163.333
$ws.Range("L140").Value = 6978.8181
$ws.Range("M140").Value = -2983.333
$ws.Range("N140").Value = -17338.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H102").Value = 1190.4445
$ws.Range("I102").Value = 1032.8572
$ws.Range("J102").Value = 1742
$ws.Range("K102").Value = 1032.8572
$ws.Range("L102").Value = 1742
$ws.Range("M102").Value = 589.1428000000001
$ws.Range("N102").Value = -4986
$ws.Range("H122").Value = 1382.1333
$ws.Range("I122").Value = 1214.5
$ws.Range("J122").Value = 1573.7142
$ws.Range("K122").Value = 3643.5
$ws.Range("L122").Value = 4721.142599999999
$ws.Range("M122").Value = -1193.5
$ws.Range("N122").Value = -9621.142599999999
$ws.Range("H132").Value = 71125.27
$ws.Range("I132").Value = 60589.47
$ws.Range("J132").Value = 86051
$ws.Range("K132").Value = 181768.41
$ws.Range("L132").Value = 258153
$ws.Range("M132").Value = -179238.41
$ws.Range("N132").Value = -263213

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 27481.514
$ws.Range("I132").Value = 1376.6786
$ws.Range("J132").Value = 93930.17999999999
$ws.Range("K132").Value = 4130.0358
$ws.Range("L132").Value = 281790.54
$ws.Range("M132").Value = -1600.0358
$ws.Range("N132").Value = -286850.54

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 59649.37
$ws.Range("I136").Value = 46296.273
$ws.Range("J136").Value = 82246.92
$ws.Range("K136").Value = 138888.819
$ws.Range("L136").Value = 246740.76
$ws.Range("M136").Value = -136338.819
$ws.Range("N136").Value = -251840.76
